$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the format (style) from the existing last header cell (O1) onto the two
# newly introduced header cells P1/Q1, without creating new style entries.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set the new header values in P1 and Q1 (row 1)
$ws.Range("P1").Value2 = 14
$ws.Range("Q1").Value2 = 15

# For each data row (2..25), set new P and Q columns to 2, and swap I/K/M/O values
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 16).Value2 = 2  # P
    $ws.Cells.Item($r, 17).Value2 = 2  # Q

    # Swap I (col 9) and K (col 11)
    $iVal = $ws.Cells.Item($r, 9).Value2
    $kVal = $ws.Cells.Item($r, 11).Value2
    $ws.Cells.Item($r, 9).Value2 = $kVal
    $ws.Cells.Item($r, 11).Value2 = $iVal

    # Swap M (col 13) and O (col 15)
    $mVal = $ws.Cells.Item($r, 13).Value2
    $oVal = $ws.Cells.Item($r, 15).Value2
    $ws.Cells.Item($r, 13).Value2 = $oVal
    $ws.Cells.Item($r, 15).Value2 = $mVal
}
